$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 <- values previously in Row 13
$ws.Cells.Item(3, 4).Value = 44819   # D3 Fecha
$ws.Cells.Item(3, 10).Value = 70  # J3 Volumen
$ws.Cells.Item(3, 11).Value = 22000  # K3 Precio minimo
$ws.Cells.Item(3, 12).Value = 22000  # L3 Precio maximo
$ws.Cells.Item(3, 13).Value = 22000  # M3 Precio promedio ponderado
$ws.Cells.Item(3, 16).Value = 1467  # P3 Precio $/Kg

# Row 4 <- values previously in Row 16
$ws.Cells.Item(4, 4).Value = 44757   # D4 Fecha
$ws.Cells.Item(4, 10).Value = 80  # J4 Volumen
$ws.Cells.Item(4, 11).Value = 25000  # K4 Precio minimo
$ws.Cells.Item(4, 12).Value = 25000  # L4 Precio maximo
$ws.Cells.Item(4, 13).Value = 25000  # M4 Precio promedio ponderado
$ws.Cells.Item(4, 16).Value = 1667  # P4 Precio $/Kg

# Row 5 <- values previously in Row 24
$ws.Cells.Item(5, 4).Value = 44418   # D5 Fecha
$ws.Cells.Item(5, 10).Value = 90  # J5 Volumen
$ws.Cells.Item(5, 11).Value = 25000  # K5 Precio minimo
$ws.Cells.Item(5, 12).Value = 25000  # L5 Precio maximo
$ws.Cells.Item(5, 13).Value = 25000  # M5 Precio promedio ponderado
$ws.Cells.Item(5, 16).Value = 1667  # P5 Precio $/Kg

# Row 6 <- values previously in Row 18
$ws.Cells.Item(6, 4).Value = 44838   # D6 Fecha
$ws.Cells.Item(6, 10).Value = 80  # J6 Volumen
$ws.Cells.Item(6, 11).Value = 22000  # K6 Precio minimo
$ws.Cells.Item(6, 12).Value = 22000  # L6 Precio maximo
$ws.Cells.Item(6, 13).Value = 22000  # M6 Precio promedio ponderado
$ws.Cells.Item(6, 16).Value = 1467  # P6 Precio $/Kg

# Row 7 <- values previously in Row 26
$ws.Cells.Item(7, 4).Value = 44778   # D7 Fecha
$ws.Cells.Item(7, 10).Value = 120  # J7 Volumen
$ws.Cells.Item(7, 11).Value = 24000  # K7 Precio minimo
$ws.Cells.Item(7, 12).Value = 24000  # L7 Precio maximo
$ws.Cells.Item(7, 13).Value = 24000  # M7 Precio promedio ponderado
$ws.Cells.Item(7, 16).Value = 1600  # P7 Precio $/Kg

# Row 8 <- values previously in Row 11
$ws.Cells.Item(8, 4).Value = 44761   # D8 Fecha
$ws.Cells.Item(8, 10).Value = 100  # J8 Volumen
$ws.Cells.Item(8, 11).Value = 23000  # K8 Precio minimo
$ws.Cells.Item(8, 12).Value = 25000  # L8 Precio maximo
$ws.Cells.Item(8, 13).Value = 24000  # M8 Precio promedio ponderado
$ws.Cells.Item(8, 16).Value = 1600  # P8 Precio $/Kg

# Row 9 <- values previously in Row 20
$ws.Cells.Item(9, 4).Value = 44400   # D9 Fecha
$ws.Cells.Item(9, 10).Value = 80  # J9 Volumen
$ws.Cells.Item(9, 11).Value = 25000  # K9 Precio minimo
$ws.Cells.Item(9, 12).Value = 25000  # L9 Precio maximo
$ws.Cells.Item(9, 13).Value = 25000  # M9 Precio promedio ponderado
$ws.Cells.Item(9, 16).Value = 1667  # P9 Precio $/Kg

# Row 11 <- values previously in Row 6
$ws.Cells.Item(11, 4).Value = 44750   # D11 Fecha
$ws.Cells.Item(11, 10).Value = 90  # J11 Volumen
$ws.Cells.Item(11, 11).Value = 25000  # K11 Precio minimo
$ws.Cells.Item(11, 12).Value = 25000  # L11 Precio maximo
$ws.Cells.Item(11, 13).Value = 25000  # M11 Precio promedio ponderado
$ws.Cells.Item(11, 16).Value = 1667  # P11 Precio $/Kg

# Row 12 <- values previously in Row 21
$ws.Cells.Item(12, 4).Value = 44365   # D12 Fecha
$ws.Cells.Item(12, 10).Value = 80  # J12 Volumen
$ws.Cells.Item(12, 11).Value = 25000  # K12 Precio minimo
$ws.Cells.Item(12, 12).Value = 25000  # L12 Precio maximo
$ws.Cells.Item(12, 13).Value = 25000  # M12 Precio promedio ponderado
$ws.Cells.Item(12, 16).Value = 1667  # P12 Precio $/Kg

# Row 13 <- values previously in Row 19
$ws.Cells.Item(13, 4).Value = 44782   # D13 Fecha
$ws.Cells.Item(13, 10).Value = 120  # J13 Volumen
$ws.Cells.Item(13, 11).Value = 24000  # K13 Precio minimo
$ws.Cells.Item(13, 12).Value = 24000  # L13 Precio maximo
$ws.Cells.Item(13, 13).Value = 24000  # M13 Precio promedio ponderado
$ws.Cells.Item(13, 16).Value = 1600  # P13 Precio $/Kg

# Row 14 <- values previously in Row 22
$ws.Cells.Item(14, 4).Value = 44754   # D14 Fecha
$ws.Cells.Item(14, 10).Value = 90  # J14 Volumen
$ws.Cells.Item(14, 11).Value = 25000  # K14 Precio minimo
$ws.Cells.Item(14, 12).Value = 25000  # L14 Precio maximo
$ws.Cells.Item(14, 13).Value = 25000  # M14 Precio promedio ponderado
$ws.Cells.Item(14, 16).Value = 1667  # P14 Precio $/Kg

# Row 15 <- values previously in Row 5
$ws.Cells.Item(15, 4).Value = 44803   # D15 Fecha
$ws.Cells.Item(15, 10).Value = 90  # J15 Volumen
$ws.Cells.Item(15, 11).Value = 24000  # K15 Precio minimo
$ws.Cells.Item(15, 12).Value = 24000  # L15 Precio maximo
$ws.Cells.Item(15, 13).Value = 24000  # M15 Precio promedio ponderado
$ws.Cells.Item(15, 16).Value = 1600  # P15 Precio $/Kg

# Row 16 <- values previously in Row 27
$ws.Cells.Item(16, 4).Value = 44740   # D16 Fecha
$ws.Cells.Item(16, 10).Value = 90  # J16 Volumen
$ws.Cells.Item(16, 11).Value = 25000  # K16 Precio minimo
$ws.Cells.Item(16, 12).Value = 25000  # L16 Precio maximo
$ws.Cells.Item(16, 13).Value = 25000  # M16 Precio promedio ponderado
$ws.Cells.Item(16, 16).Value = 1667  # P16 Precio $/Kg

# Row 17 <- values previously in Row 15
$ws.Cells.Item(17, 4).Value = 44810   # D17 Fecha
$ws.Cells.Item(17, 10).Value = 110  # J17 Volumen
$ws.Cells.Item(17, 11).Value = 22000  # K17 Precio minimo
$ws.Cells.Item(17, 12).Value = 22000  # L17 Precio maximo
$ws.Cells.Item(17, 13).Value = 22000  # M17 Precio promedio ponderado
$ws.Cells.Item(17, 16).Value = 1467  # P17 Precio $/Kg

# Row 18 <- values previously in Row 14
$ws.Cells.Item(18, 4).Value = 44781   # D18 Fecha
$ws.Cells.Item(18, 10).Value = 70  # J18 Volumen
$ws.Cells.Item(18, 11).Value = 24000  # K18 Precio minimo
$ws.Cells.Item(18, 12).Value = 24000  # L18 Precio maximo
$ws.Cells.Item(18, 13).Value = 24000  # M18 Precio promedio ponderado
$ws.Cells.Item(18, 16).Value = 1600  # P18 Precio $/Kg

# Row 19 <- values previously in Row 8
$ws.Cells.Item(19, 4).Value = 44789   # D19 Fecha
$ws.Cells.Item(19, 10).Value = 90  # J19 Volumen
$ws.Cells.Item(19, 11).Value = 24000  # K19 Precio minimo
$ws.Cells.Item(19, 12).Value = 24000  # L19 Precio maximo
$ws.Cells.Item(19, 13).Value = 24000  # M19 Precio promedio ponderado
$ws.Cells.Item(19, 16).Value = 1600  # P19 Precio $/Kg

# Row 20 <- values previously in Row 17
$ws.Cells.Item(20, 4).Value = 44775   # D20 Fecha
$ws.Cells.Item(20, 10).Value = 120  # J20 Volumen
$ws.Cells.Item(20, 11).Value = 24000  # K20 Precio minimo
$ws.Cells.Item(20, 12).Value = 24000  # L20 Precio maximo
$ws.Cells.Item(20, 13).Value = 24000  # M20 Precio promedio ponderado
$ws.Cells.Item(20, 16).Value = 1600  # P20 Precio $/Kg

# Row 21 <- values previously in Row 3
$ws.Cells.Item(21, 4).Value = 44817   # D21 Fecha
$ws.Cells.Item(21, 10).Value = 90  # J21 Volumen
$ws.Cells.Item(21, 11).Value = 23000  # K21 Precio minimo
$ws.Cells.Item(21, 12).Value = 23000  # L21 Precio maximo
$ws.Cells.Item(21, 13).Value = 23000  # M21 Precio promedio ponderado
$ws.Cells.Item(21, 16).Value = 1533  # P21 Precio $/Kg

# Row 22 <- values previously in Row 9
$ws.Cells.Item(22, 4).Value = 44831   # D22 Fecha
$ws.Cells.Item(22, 10).Value = 90  # J22 Volumen
$ws.Cells.Item(22, 11).Value = 25000  # K22 Precio minimo
$ws.Cells.Item(22, 12).Value = 25000  # L22 Precio maximo
$ws.Cells.Item(22, 13).Value = 25000  # M22 Precio promedio ponderado
$ws.Cells.Item(22, 16).Value = 1667  # P22 Precio $/Kg

# Row 24 <- values previously in Row 12
$ws.Cells.Item(24, 4).Value = 44799   # D24 Fecha
$ws.Cells.Item(24, 10).Value = 80  # J24 Volumen
$ws.Cells.Item(24, 11).Value = 23000  # K24 Precio minimo
$ws.Cells.Item(24, 12).Value = 23000  # L24 Precio maximo
$ws.Cells.Item(24, 13).Value = 23000  # M24 Precio promedio ponderado
$ws.Cells.Item(24, 16).Value = 1533  # P24 Precio $/Kg

# Row 25 <- values previously in Row 4
$ws.Cells.Item(25, 4).Value = 44407   # D25 Fecha
$ws.Cells.Item(25, 10).Value = 90  # J25 Volumen
$ws.Cells.Item(25, 11).Value = 25000  # K25 Precio minimo
$ws.Cells.Item(25, 12).Value = 25000  # L25 Precio maximo
$ws.Cells.Item(25, 13).Value = 25000  # M25 Precio promedio ponderado
$ws.Cells.Item(25, 16).Value = 1667  # P25 Precio $/Kg

# Row 26 <- values previously in Row 25
$ws.Cells.Item(26, 4).Value = 44764   # D26 Fecha
$ws.Cells.Item(26, 10).Value = 90  # J26 Volumen
$ws.Cells.Item(26, 11).Value = 24000  # K26 Precio minimo
$ws.Cells.Item(26, 12).Value = 24000  # L26 Precio maximo
$ws.Cells.Item(26, 13).Value = 24000  # M26 Precio promedio ponderado
$ws.Cells.Item(26, 16).Value = 1600  # P26 Precio $/Kg

# Row 27 <- values previously in Row 7
$ws.Cells.Item(27, 4).Value = 44792   # D27 Fecha
$ws.Cells.Item(27, 10).Value = 120  # J27 Volumen
$ws.Cells.Item(27, 11).Value = 24000  # K27 Precio minimo
$ws.Cells.Item(27, 12).Value = 24000  # L27 Precio maximo
$ws.Cells.Item(27, 13).Value = 24000  # M27 Precio promedio ponderado
$ws.Cells.Item(27, 16).Value = 1600  # P27 Precio $/Kg

